$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H76").Value = 5103.769
$ws.Range("I76").Value = 3950
$ws.Range("J76").Value = 5616.5557
$ws.Range("K76").Value = 3950
$ws.Range("L76").Value = 5616.5557
$ws.Range("M76").Value = -3635
$ws.Range("N76").Value = -6246.5557

$ws.Range("H79").Value = 5103.769
$ws.Range("I79").Value = 3950
$ws.Range("J79").Value = 5616.5557
$ws.Range("K79").Value = 3950
$ws.Range("L79").Value = 5616.5557
$ws.Range("M79").Value = -2858
$ws.Range("N79").Value = -7800.5557

$ws.Range("H80").Value = 860.7895
$ws.Range("I80").Value = 612.6667
$ws.Range("J80").Value = 1286.1428
$ws.Range("K80").Value = 1838.0001
$ws.Range("L80").Value = 3858.4284
$ws.Range("M80").Value = -840.0001
$ws.Range("N80").Value = -5854.428400000001

$ws.Range("H83").Value = 860.7895
$ws.Range("I83").Value = 612.6667
$ws.Range("J83").Value = 1286.1428
$ws.Range("K83").Value = 5514.0003
$ws.Range("L83").Value = 11575.2852
$ws.Range("M83").Value = -522.0002999999997
$ws.Range("N83").Value = -21559.2852

$ws.Range("H92").Value = 747.2
$ws.Range("I92").Value = 396
$ws.Range("J92").Value = 1274
$ws.Range("K92").Value = 396
$ws.Range("L92").Value = 1274
$ws.Range("M92").Value = 852
$ws.Range("N92").Value = -3770

$ws.Range("H137").Value = 3847.1333
$ws.Range("I137").Value = 1138.0625
$ws.Range("J137").Value = 5341.793
$ws.Range("K137").Value = 3414.1875
$ws.Range("L137").Value = 16025.379
$ws.Range("M137").Value = -864.1875
$ws.Range("N137").Value = -21125.379

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H88").Value = 2191.2222
$ws.Range("I88").Value = 2204.8
$ws.Range("J88").Value = 2174.25
$ws.Range("K88").Value = 2204.8
$ws.Range("L88").Value = 2174.25
$ws.Range("M88").Value = -1798.8
$ws.Range("N88").Value = -2986.25

$ws.Range("H91").Value = 2191.2222
$ws.Range("I91").Value = 2204.8
$ws.Range("J91").Value = 2174.25
$ws.Range("K91").Value = 2204.8
$ws.Range("L91").Value = 2174.25
$ws.Range("M91").Value = -800.8000000000002
$ws.Range("N91").Value = -4982.25

$ws.Range("H132").Value = 15193.9375
$ws.Range("I132").Value = 16851.857
$ws.Range("J132").Value = 13904.444
$ws.Range("K132").Value = 50555.571
$ws.Range("L132").Value = 41713.33199999999
$ws.Range("M132").Value = -48025.571
$ws.Range("N132").Value = -46773.33199999999

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H8").Value = 0
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 0
$ws.Range("M8").ClearContents()

$ws.Range("H20").Value = 63390.875
$ws.Range("I20").Value = 995
$ws.Range("J20").Value = 200661.8
$ws.Range("K20").Value = 995
$ws.Range("L20").Value = 200661.8
$ws.Range("M20").Value = -748
$ws.Range("N20").Value = -201155.8

$ws.Range("H86").Value = 1782.5333
$ws.Range("I86").Value = 1313.4546
$ws.Range("J86").Value = 3072.5
$ws.Range("K86").Value = 1313.4546
$ws.Range("L86").Value = 3072.5
$ws.Range("M86").Value = -190.4546
$ws.Range("N86").Value = -5318.5

$ws.Range("H89").Value = 1782.5333
$ws.Range("I89").Value = 1313.4546
$ws.Range("J89").Value = 3072.5
$ws.Range("K89").Value = 6567.273
$ws.Range("L89").Value = 15362.5
$ws.Range("M89").Value = -951.2730000000001
$ws.Range("N89").Value = -26594.5

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H2").Value = 0
$ws.Range("I2").Value = 0
$ws.Range("K2").Value = 0
$ws.Range("M2").ClearContents()

$ws.Range("H31").Value = 30007.4
$ws.Range("I31").Value = 1846.5264
$ws.Range("J31").Value = 55486.285
$ws.Range("K31").Value = 1846.5264
$ws.Range("L31").Value = 55486.285
$ws.Range("M31").Value = -1551.5264
$ws.Range("N31").Value = -56076.285

$ws.Range("H34").Value = 30007.4
$ws.Range("I34").Value = 1846.5264
$ws.Range("J34").Value = 55486.285
$ws.Range("K34").Value = 1846.5264
$ws.Range("L34").Value = 55486.285
$ws.Range("M34").Value = -1644.5264
$ws.Range("N34").Value = -55890.285

$ws.Range("H62").Value = 4177.75
$ws.Range("J62").Value = 2699.6667
$ws.Range("L62").Value = 2699.6667
$ws.Range("N62").Value = -3947.6667

$ws.Range("H65").Value = 4177.75
$ws.Range("J65").Value = 2699.6667
$ws.Range("L65").Value = 13498.3335
$ws.Range("N65").Value = -19738.3335

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H82").Value = 1900
$ws.Range("I82").Value = 1900
$ws.Range("J82").Value = 0
$ws.Range("K82").Value = 5700
$ws.Range("L82").Value = 0
$ws.Range("M82").Value = -5294
$ws.Range("N82").ClearContents()

$ws.Range("H85").Value = 1900
$ws.Range("I85").Value = 1900
$ws.Range("J85").Value = 0
$ws.Range("K85").Value = 5700
$ws.Range("L85").Value = 0
$ws.Range("M85").Value = -4296
$ws.Range("N85").ClearContents()

$ws.Range("H92").Value = 929
$ws.Range("J92").Value = 1017.1667
$ws.Range("L92").Value = 3051.5001
$ws.Range("N92").Value = -5547.5001

$ws.Range("H93").Value = 7888
$ws.Range("I93").Value = 1552
$ws.Range("J93").Value = 10000
$ws.Range("K93").Value = 4656
$ws.Range("L93").Value = 30000
$ws.Range("M93").Value = -2784
$ws.Range("N93").Value = -33744

$ws.Range("H94").Value = 0
$ws.Range("I94").Value = 0
$ws.Range("K94").Value = 0
$ws.Range("M94").ClearContents()

$ws.Range("H101").Value = 7750
$ws.Range("J101").Value = 7750
$ws.Range("L101").Value = 23250
$ws.Range("N101").Value = -28118

$ws.Range("H102").Value = 3126.8333
$ws.Range("I102").Value = 130
$ws.Range("J102").Value = 3399.2727
$ws.Range("K102").Value = 390
$ws.Range("L102").Value = 10197.8181
$ws.Range("M102").Value = 2044
$ws.Range("N102").Value = -15065.8181

$ws.Range("H103").Value = 80
$ws.Range("I103").Value = 80
$ws.Range("J103").Value = 0
$ws.Range("K103").Value = 240
$ws.Range("L103").Value = 0
$ws.Range("M103").Value = 639
$ws.Range("N103").ClearContents()

$ws.Range("H105").Value = 12547.692
$ws.Range("J105").Value = 12547.692
$ws.Range("L105").Value = 37643.076
$ws.Range("N105").Value = -42885.076

$ws.Range("H106").Value = 2862.5
$ws.Range("J106").Value = 2862.5
$ws.Range("L106").Value = 8587.5
$ws.Range("N106").Value = -10479.5

$ws.Range("H107").Value = 308
$ws.Range("J107").Value = 358.5
$ws.Range("L107").Value = 1075.5
$ws.Range("N107").Value = -4915.5

$ws.Range("H108").Value = 589.75
$ws.Range("I108").Value = 245.42857
$ws.Range("J108").Value = 3000
$ws.Range("K108").Value = 736.28571
$ws.Range("L108").Value = 9000
$ws.Range("M108").Value = 2143.71429
$ws.Range("N108").Value = -14760

$ws.Range("H109").Value = 0
$ws.Range("I109").Value = 0
$ws.Range("K109").Value = 0
$ws.Range("M109").ClearContents()

$ws.Range("H113").Value = 506.2
$ws.Range("I113").Value = 0
$ws.Range("J113").Value = 506.2
$ws.Range("K113").Value = 0
$ws.Range("L113").Value = 1518.6
$ws.Range("N113").Value = -5858.6
$ws.Range("M113").ClearContents()

$ws.Range("H114").Value = 0
$ws.Range("I114").Value = 0
$ws.Range("K114").Value = 0
$ws.Range("M114").ClearContents()

$ws.Range("H115").Value = 1400
$ws.Range("I115").Value = 1400
$ws.Range("K115").Value = 4200
$ws.Range("M115").Value = -3025

$ws.Range("H116").Value = 3257.25
$ws.Range("I116").Value = 3676.3333
$ws.Range("K116").Value = 11028.9999
$ws.Range("M116").Value = -7586.999899999999

$ws.Range("H117").Value = 385.66666
$ws.Range("I117").Value = 385.66666
$ws.Range("J117").Value = 0
$ws.Range("K117").Value = 1156.99998
$ws.Range("L117").Value = 0
$ws.Range("M117").Value = 2285.00002
$ws.Range("N117").ClearContents()

$ws.Range("H118").Value = 3125
$ws.Range("I118").Value = 0
$ws.Range("J118").Value = 3125
$ws.Range("K118").Value = 0
$ws.Range("L118").Value = 9375
$ws.Range("N118").Value = -11861
$ws.Range("M118").ClearContents()

$ws.Range("H121").Value = 12580.333
$ws.Range("I121").Value = 0
$ws.Range("J121").Value = 12580.333
$ws.Range("K121").Value = 0
$ws.Range("L121").Value = 37740.999
$ws.Range("N121").Value = -40360.999
$ws.Range("M121").ClearContents()

$ws.Range("H122").Value = 1347.4445
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 1347.4445
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 12127.0005
$ws.Range("N122").Value = -17027.0005
$ws.Range("M122").ClearContents()

$ws.Range("H126").Value = 0
$ws.Range("I126").Value = 0
$ws.Range("J126").Value = 0
$ws.Range("K126").Value = 0
$ws.Range("L126").Value = 0
$ws.Range("M126").ClearContents()
$ws.Range("N126").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H70").Value = 4551.36
$ws.Range("I70").Value = 4476.143
$ws.Range("J70").Value = 4647.091
$ws.Range("K70").Value = 4476.143
$ws.Range("L70").Value = 4647.091
$ws.Range("M70").Value = -4206.143
$ws.Range("N70").Value = -5187.091

$ws.Range("H73").Value = 4551.36
$ws.Range("I73").Value = 4476.143
$ws.Range("J73").Value = 4647.091
$ws.Range("K73").Value = 4476.143
$ws.Range("L73").Value = 4647.091
$ws.Range("M73").Value = -3540.143
$ws.Range("N73").Value = -6519.091

$ws.Range("H80").Value = 3261.1333
$ws.Range("I80").Value = 2205
$ws.Range("J80").Value = 3525.1667
$ws.Range("K80").Value = 2205
$ws.Range("L80").Value = 3525.1667
$ws.Range("M80").Value = -1207
$ws.Range("N80").Value = -5521.1667

$ws.Range("H83").Value = 3261.1333
$ws.Range("I83").Value = 2205
$ws.Range("J83").Value = 3525.1667
$ws.Range("K83").Value = 11025
$ws.Range("L83").Value = 17625.8335
$ws.Range("M83").Value = -6033
$ws.Range("N83").Value = -27609.8335

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H5").Value = 5755.385
$ws.Range("I5").Value = 2970
$ws.Range("J5").Value = 5987.5
$ws.Range("K5").Value = 2970
$ws.Range("L5").Value = 5987.5
$ws.Range("M5").Value = -2858
$ws.Range("N5").Value = -6211.5

$ws.Range("H113").Value = 43884.824
$ws.Range("I113").Value = 62832.688
$ws.Range("J113").Value = 575.4286
$ws.Range("K113").Value = 188498.064
$ws.Range("L113").Value = 1726.2858
$ws.Range("M113").Value = -186328.064
$ws.Range("N113").Value = -6066.2858
